$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Fill in the previously-empty N/O (interregional) columns on sheet 1 (mse (unbalanced set))
$ws1.Range("N6").Value = 0.0047695003077224078
$ws1.Range("O6").Value = 0.03118292825736188
$ws1.Range("N7").Value = 0.0045720456292513336
$ws1.Range("O7").Value = 0.033525868579185748
$ws1.Range("N8").Value = 0.0045515442489738557
$ws1.Range("O8").Value = 0.033755682566183873
$ws1.Range("N9").Value = 0.0046259677353847383
$ws1.Range("O9").Value = 0.033223096939218863
$ws1.Range("N10").Value = 0.0045361807639645356
$ws1.Range("O10").Value = 0.031339850159775147
$ws1.Range("N11").Value = 0.0045716462394260133
$ws1.Range("O11").Value = 0.035131207937474568
$ws1.Range("N12").Value = 0.0045960000125573894
$ws1.Range("O12").Value = 0.033711821218537587
$ws1.Range("N13").Value = 0.0046727769605473214
$ws1.Range("O13").Value = 0.033532651304816641
$ws1.Range("N14").Value = 0.0045249634270531636
$ws1.Range("O14").Value = 0.038052163397268322
$ws1.Range("N15").Value = 0.0046851858352022907
$ws1.Range("O15").Value = 0.0295066706798136
$ws1.Range("N16").Value = 0.0046711126068756043
$ws1.Range("O16").Value = 0.031890609473540422
$ws1.Range("N17").Value = 0.004595814503362672
$ws1.Range("O17").Value = 0.031479945387102873
$ws1.Range("N18").Value = 0.0045985444195138928
$ws1.Range("O18").Value = 0.032261766252048477
$ws1.Range("N19").Value = 0.0046672749040183176
$ws1.Range("O19").Value = 0.031205521975054051
$ws1.Range("N20").Value = 0.0046756312767027566
$ws1.Range("O20").Value = 0.030543567477029351
$ws1.Range("N21").Value = 0.0047327213927244594
$ws1.Range("O21").Value = 0.034007994157639679
$ws1.Range("N22").Value = 0.0046599265795111587
$ws1.Range("O22").Value = 0.032751649259127082
$ws1.Range("N23").Value = 0.0044404871608199259
$ws1.Range("O23").Value = 0.037258792699951603
$ws1.Range("N24").Value = 0.0046850595714198551
$ws1.Range("O24").Value = 0.033063752032043793
$ws1.Range("N25").Value = 0.0045944024860740213
$ws1.Range("O25").Value = 0.03186562841467093
$ws1.Range("N26").Value = 0.0047206878161942306
$ws1.Range("O26").Value = 0.030976917866293009
$ws1.Range("N27").Value = 0.0046768410862040814
$ws1.Range("O27").Value = 0.031723597198414458
$ws1.Range("N28").Value = 0.0046595010932994826
$ws1.Range("O28").Value = 0.031633979506510833
$ws1.Range("N29").Value = 0.0046488932782003346
$ws1.Range("O29").Value = 0.035018188400716688
$ws1.Range("N30").Value = 0.0044754589531873293
$ws1.Range("O30").Value = 0.035571871154928172
$ws1.Range("N31").Value = 0.0048468629148511293
$ws1.Range("O31").Value = 0.029428253484147342
$ws1.Range("N32").Value = 0.0045299699096176338
$ws1.Range("O32").Value = 0.033993385068867177
$ws1.Range("N33").Value = 0.0044468747734298332
$ws1.Range("O33").Value = 0.035616620067012897
$ws1.Range("N34").Value = 0.0046087905955834573
$ws1.Range("O34").Value = 0.03427787405639296
$ws1.Range("N35").Value = 0.004560760129113382
$ws1.Range("O35").Value = 0.031577120832329247
$ws1.Range("N36").Value = 0.0045141551074862566
$ws1.Range("O36").Value = 0.03568596345990286
$ws1.Range("N37").Value = 0.0045930025186855317
$ws1.Range("O37").Value = 0.033605636319693898
$ws1.Range("N38").Value = 0.0045747805887478774
$ws1.Range("O38").Value = 0.029587805335058889
$ws1.Range("N39").Value = 0.0045343388878611914
$ws1.Range("O39").Value = 0.034012705798304579
$ws1.Range("N40").Value = 0.0046604170580406574
$ws1.Range("O40").Value = 0.031516542249970768
$ws1.Range("N41").Value = 0.004623230705287292
$ws1.Range("O41").Value = 0.03118880653929898
$ws1.Range("N42").Value = 0.0046227011351756987
$ws1.Range("O42").Value = 0.035455802631071648
$ws1.Range("N43").Value = 0.0046537650321949384
$ws1.Range("O43").Value = 0.032489236391748243
$ws1.Range("N44").Value = 0.0046426568961247224
$ws1.Range("O44").Value = 0.03341286040584953
$ws1.Range("N45").Value = 0.0044332302146524944
$ws1.Range("O45").Value = 0.038752375060404833
$ws1.Range("N46").Value = 0.004673619534517324
$ws1.Range("O46").Value = 0.02936176879549791
$ws1.Range("N47").Value = 0.0045802204049157176
$ws1.Range("O47").Value = 0.032698903014279683
$ws1.Range("N48").Value = 0.0047043106275732986
$ws1.Range("O48").Value = 0.030282751698742449
$ws1.Range("N49").Value = 0.0046604296739450992
$ws1.Range("O49").Value = 0.032144299105788747
$ws1.Range("N50").Value = 0.0045346488444017114
$ws1.Range("O50").Value = 0.036354363943475479
$ws1.Range("N51").Value = 0.0045783474977289414
$ws1.Range("O51").Value = 0.031746551808347588
$ws1.Range("N52").Value = 0.0046121546506134106
$ws1.Range("O52").Value = 0.033317772709260049
$ws1.Range("N53").Value = 0.0044702526660844004
$ws1.Range("O53").Value = 0.035742529622319587
$ws1.Range("N54").Value = 0.004694499494012701
$ws1.Range("O54").Value = 0.031921624442287702
$ws1.Range("N55").Value = 0.0044949919681506788
$ws1.Range("O55").Value = 0.03359791864684826

# Prepare the new F column (interregional) on sheet 2 (feature sig) with the same
# number format / alignment used by the neighbouring D and E columns, then fill values
$ws2.Range("F5:F19").NumberFormat = $ws2.Range("D5").NumberFormat
$ws2.Range("F5:F19").HorizontalAlignment = $ws2.Range("D5").HorizontalAlignment
$ws2.Range("F5").Value = 0.088413027272881861
$ws2.Range("F6").Value = 0.13328645042113921
$ws2.Range("F7").Value = 0.079719945271623127
$ws2.Range("F8").Value = 0.06058922124692831
$ws2.Range("F9").Value = 0.056071573053039482
$ws2.Range("F10").Value = 0.066070805140619762
$ws2.Range("F11").Value = 0.055571882969546488
$ws2.Range("F12").Value = 0.053242116681780208
$ws2.Range("F13").Value = 0.053475660887844681
$ws2.Range("F14").Value = 0.057848632447795943
$ws2.Range("F15").Value = 0.06580257898426678
$ws2.Range("F16").Value = 0.060617071528018117
$ws2.Range("F17").Value = 0.066874939593002225
$ws2.Range("F18").Value = 0.037055712634370012
$ws2.Range("F19").Value = 0.065360381867143935

# Add page setup (paper size / orientation) to sheet 2
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Update selections: sheet2 first so sheet1 ends up the active/tabSelected sheet
$ws2.Range("G19").Select()
$ws1.Range("V43").Select()
